$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with columns M:V, reusing the existing
#     bold/bordered header style from B1 (same style as C1:E1). ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("M1:V1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("M1").Value = 11
$ws.Range("N1").Value = 12
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16
$ws.Range("S1").Value = 17
$ws.Range("T1").Value = 18
$ws.Range("U1").Value = 19
$ws.Range("V1").Value = 20

# --- Row 2: Final Value results for the new models ---
$ws.Range("M2").Value = 2077777.6039563499
$ws.Range("N2").Value = 1299889.089827256
$ws.Range("O2").Value = 1831931.490443537
$ws.Range("P2").Value = 2398477.3009716091
$ws.Range("Q2").Value = 2457431.301863702
$ws.Range("R2").Value = 935888.04287778307
$ws.Range("S2").Value = 2341385.2750639911
$ws.Range("T2").Value = 2611651.4317320092
$ws.Range("U2").Value = 1738344.1965679959
$ws.Range("V2").Value = 3043745.4990403811

# --- Row 3: Annualized Return results for the new models ---
$ws.Range("M3").Value = 0.27319552246481038
$ws.Range("N3").Value = 0.090486840098330523
$ws.Range("O3").Value = 0.2213282582508154
$ws.Range("P3").Value = 0.33500631981547252
$ws.Range("Q3").Value = 0.34575601393580541
$ws.Range("R3").Value = -0.021646131625614681
$ws.Range("S3").Value = 0.3244261318390993
$ws.Range("T3").Value = 0.37308292809971322
$ws.Range("U3").Value = 0.2003582928149894
$ws.Range("V3").Value = 0.4443012900483243

# --- Row 4: Sharpe Ratio results for the new models ---
$ws.Range("M4").Value = 0.77926087901354657
$ws.Range("N4").Value = 0.22310295603550381
$ws.Range("O4").Value = 0.41262871607050061
$ws.Range("P4").Value = 0.97773526377681963
$ws.Range("Q4").Value = 0.89282703341417369
$ws.Range("R4").Value = -0.1300086352744057
$ws.Range("S4").Value = 0.98590854898934233
$ws.Range("T4").Value = 1.213892890441558
$ws.Range("U4").Value = 0.5279717441629721
$ws.Range("V4").Value = 1.130510503050099

# --- Column W was bumped to a fixed width (11 chars) by the author,
#     even though it holds no data. ---
$ws.Columns("W").ColumnWidth = 10.1

# --- Restore the selection to where the author left the cursor ---
$ws.Range("R8").Select() | Out-Null
